# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Femacal de La Calera" / Frutilla
# right before the existing block that starts at row 637, shifting the
# existing rows (637-652) down to (641-656).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 637 (pushes old 637.. down to 641..)
$ws.Rows("637:640").Insert()

# Common (unchanged) metadata columns shared by every row in this block
$marketId   = 3
$market     = "Femacal de La Calera"
$region     = "Coquimbo"
$codreg     = 5
$tipo       = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria  = "Frutilla"
$variedad   = "Sin especificar"
$unidad     = "$/bandeja 7 kilos"
$origen     = "Provincia de Melipilla"
$kgUnidad   = 7

# New rows data: row, date(serial), calidad, volumen, precioMin, precioMax, precioProm, precioKg
$newRows = @(
    @{ Row = 637; Fecha = 45239; Calidad = "Especial"; Volumen = 143; PMin = 12000; PMax = 13000; PProm = 12524; PKg = 1789 },
    @{ Row = 638; Fecha = 45239; Calidad = "Primera";  Volumen = 68;  PMin = 10000; PMax = 10000; PProm = 10000; PKg = 1429 },
    @{ Row = 639; Fecha = 45239; Calidad = "Segunda";  Volumen = 68;  PMin = 8000;  PMax = 8000;  PProm = 8000;  PKg = 1143 },
    @{ Row = 640; Fecha = 45239; Calidad = "Tercera";  Volumen = 50;  PMin = 5000;  PMax = 5000;  PProm = 5000;  PKg = 714 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $marketId
    $ws.Cells.Item($row, 2).Value = $market
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
